$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110 — this shifts the existing row 110
# (and everything below it) down by one row, growing the used range
# from A1:R235 to A1:R236.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly price record.
$ws.Range("A110").Value = 8
$ws.Range("B110").Value = "Terminal La Palmera de La Serena"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 45118
$ws.Range("E110").Value = 4
$ws.Range("F110").Value = 100112044
$ws.Range("G110").Value = "Perejil"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2500
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 2750
$ws.Range("N110").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O110").Value = "Provincia del Elquí"
$ws.Range("P110").Value = 1833
$ws.Range("Q110").Value = 1.5
$ws.Range("R110").Value = "Hortaliza"
